$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Naming convention change: Txx -> T_x and NTxx -> NT_x for the patient
# sample columns in the header row.
$ws.Range("B1").Value = "T_1"
$ws.Range("C1").Value = "T_2"
$ws.Range("D1").Value = "T_3"
$ws.Range("E1").Value = "T_4"
$ws.Range("F1").Value = "T_5"
$ws.Range("G1").Value = "NT_1"
$ws.Range("H1").Value = "NT_2"
$ws.Range("I1").Value = "NT_3"
$ws.Range("J1").Value = "NT_4"
$ws.Range("K1").Value = "NT_5"

# The "CpG_Array" row label is renamed to "Mean_beta-value".
$ws.Range("A22").Value = "Mean_beta-value"

# Move/record the active selection on the bottom-right pane to A23, matching
# the saved view state in the workbook.
$ws.Range("A23").Select()
